$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; this shifts the existing rows 31-129
# down to 32-130, matching the target dimension A1:R130.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record's data.
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C31").Value = "Ñuble"
$ws.Range("D31").Value = 44607
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = 100112045
$ws.Range("G31").Value = "Zapallo"
$ws.Range("H31").Value = "Camote"
$ws.Range("I31").Value = "1a (cosecha)"
$ws.Range("J31").Value = 120
$ws.Range("K31").Value = 350
$ws.Range("L31").Value = 400
$ws.Range("M31").Value = 375
$ws.Range("N31").Value = '$/kilo (volumen en unidades)'
$ws.Range("O31").Value = "Región de O'Higgins"
$ws.Range("P31").Value = 375
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
